# Actualizar cotizaciones de prueba con nuevo sistema de mensajes
# Adds 3 new quote rows (24-26) to the Cotizaciones sheet, copying the
# formatting of the last existing row and filling in the new data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Add-CotizacionRow {
    param(
        $TargetRow,
        $SourceRow,
        $FechaRegistro,
        $FechaServicio,
        $HoraServicio,
        $NombreCompleto,
        $Email,
        $Telefono,
        $TelefonoEmergencia,
        $Origen,
        $Destino,
        $NParadas,
        $Distancia,
        $Duracion,
        $Personas,
        $MarcaModelo,
        $Transmision,
        $Patente,
        $Seguro,
        $CostoBase,
        $CostoFinal,
        $CodigoDescuento,
        $Descuento
    )

    # Copy the full formatting (styles/borders/fills) from the reference row
    $srcRange = $ws.Range("A$($SourceRow):Z$($SourceRow)")
    $dstRange = $ws.Range("A$($TargetRow):Z$($TargetRow)")
    $srcRange.Copy()
    $dstRange.PasteSpecial(-4122) # xlPasteFormats

    # Also copy values so the "separator" columns (D,I,L,Q,V - normally
    # blank/empty text) match the existing blank-text pattern used across
    # the sheet.
    $srcRange.Copy()
    $dstRange.PasteSpecial(-4163) # xlPasteValues

    $ws.Cells.Item($TargetRow, 1).Value = $FechaRegistro
    $ws.Cells.Item($TargetRow, 2).Value = $FechaServicio
    $ws.Cells.Item($TargetRow, 3).Value = $HoraServicio
    $ws.Cells.Item($TargetRow, 5).Value = $NombreCompleto
    $ws.Cells.Item($TargetRow, 6).Value = $Email
    $ws.Cells.Item($TargetRow, 7).Value = $Telefono
    $ws.Cells.Item($TargetRow, 8).Value = $TelefonoEmergencia
    $ws.Cells.Item($TargetRow, 10).Value = $Origen
    $ws.Cells.Item($TargetRow, 11).Value = $Destino
    $ws.Cells.Item($TargetRow, 13).Value = $NParadas
    $ws.Cells.Item($TargetRow, 14).Value = $Distancia
    $ws.Cells.Item($TargetRow, 15).Value = $Duracion
    $ws.Cells.Item($TargetRow, 16).Value = $Personas
    $ws.Cells.Item($TargetRow, 18).Value = $MarcaModelo
    $ws.Cells.Item($TargetRow, 19).Value = $Transmision
    $ws.Cells.Item($TargetRow, 20).Value = $Patente
    $ws.Cells.Item($TargetRow, 21).Value = $Seguro
    $ws.Cells.Item($TargetRow, 23).Value = $CostoBase
    $ws.Cells.Item($TargetRow, 24).Value = $CostoFinal
    $ws.Cells.Item($TargetRow, 25).Value = $CodigoDescuento
    $ws.Cells.Item($TargetRow, 26).Value = $Descuento
}

Add-CotizacionRow 24 23 `
    "03/01/2026 16:42:49" `
    "2026-01-16" `
    "01:00" `
    "benjamin Riveros " `
    "benjamin4riveros@gmail.com" `
    "956130912" `
    "956061185" `
    "Pacul, La Florida, Región Metropolitana de Santiago 8240000, Chile" `
    "Peñalolén, Región Metropolitana de Santiago, Chile" `
    0 `
    "14.93" `
    27 `
    "2" `
    "ford " `
    "mecanico" `
    "FJDG47" `
    "si" `
    32990 `
    29691 `
    "123" `
    3299

Add-CotizacionRow 25 24 `
    "03/01/2026 16:48:19" `
    "2026-01-29" `
    "23:00" `
    "benjamin Riveros " `
    "benjamin4riveros@gmail.com" `
    "956130912" `
    "956061185" `
    "Pacul, La Florida, Región Metropolitana de Santiago 8240000, Chile" `
    "Peñalolén, Región Metropolitana de Santiago, Chile" `
    0 `
    "14.93" `
    27 `
    "2" `
    "ford " `
    "automatico" `
    "FJDG47" `
    "si" `
    32990 `
    29691 `
    "123" `
    3299

Add-CotizacionRow 26 25 `
    "03/01/2026 16:58:19" `
    "2026-01-24" `
    "05:00" `
    "benjamin Riveros " `
    "benjamin4riveros@gmail.com" `
    "956130912" `
    "956130912" `
    "Pacul, La Florida, Región Metropolitana de Santiago 8240000, Chile" `
    "Peñalolén, Región Metropolitana de Santiago, Chile" `
    0 `
    "12.94" `
    26 `
    "2" `
    "ford " `
    "automatico" `
    "FJDG47" `
    "si" `
    31990 `
    28791 `
    "123" `
    3199

Write-Host "Added rows 24-26 to Cotizaciones sheet"
